$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "70.028.05"
$ws.Range("E2").Value = "  -0.60%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.512.20"
$ws.Range("E3").Value = "  -1.63%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "575.77"
$ws.Range("E5").Value = "  -0.25%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "166.03"
$ws.Range("E6").Value = "  -2.45%  "
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.520"
$ws.Range("E8").Value = "  +2.00%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.508.04"
$ws.Range("E9").Value = "  -1.77%  "
$ws.Range("E10").Value = "  -2.34%  "
$ws.Range("E11").Value = "  -1.08%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.355"
$ws.Range("E12").Value = "  +2.75%  "
$ws.Range("E13").Value = "  +1.78%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.974.37"
$ws.Range("E14").Value = "  -1.56%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "69.934.81"
$ws.Range("E15").Value = "  -0.62%  "
$ws.Range("E16").Value = "  -2.61%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "24.99"
$ws.Range("E17").Value = "  -0.55%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.532.70"
$ws.Range("E18").Value = "  -0.41%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.41"
$ws.Range("E19").Value = "  -2.64%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.78"
$ws.Range("E20").Value = "  +1.06%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "350.93"
$ws.Range("E21").Value = "  -2.90%  "
$ws.Range("E22").Value = "  -1.16%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.98"
$ws.Range("E23").Value = "  -1.60%  "
$ws.Range("E24").Value = "  +0.10%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "70.32"
$ws.Range("E25").Value = "  +0.48%  "
$ws.Range("E26").Value = "  -2.02%  "
$ws.Range("B27").Value = "Aptos"
$ws.Range("C27").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.85"
$ws.Range("E27").Value = "  -5.09%  "
$ws.Range("B28").Value = "WrappedeETH"
$ws.Range("C28").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.657.00"
$ws.Range("E28").Value = "  -1.23%  "
$ws.Range("E29").Value = "  +0.32%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0₃0899"
$ws.Range("E30").Value = "  -3.28%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.88"
$ws.Range("E31").Value = "  -0.03%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "465.04"
$ws.Range("E32").Value = "  -4.26%  "
$ws.Range("E33").Value = "  -3.20%  "
$ws.Range("E34").Value = "  -1.62%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.999"
$ws.Range("E35").Value = "  -0.09%  "
$ws.Range("E36").Value = "  +0.56%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "155.55"
$ws.Range("E37").Value = "  -1.38%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "19.06"
$ws.Range("E38").Value = "  +1.26%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "18.63"
$ws.Range("E39").Value = "  -0.48%  "
$ws.Range("E40").Value = "  -0.03%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "4.76"
$ws.Range("E41").Value = "  -0.43%  "
$ws.Range("E42").Value = "  -0.75%  "
$ws.Range("E43").Value = "  -4.14%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "38.43"
$ws.Range("E45").Value = "  -13.37%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.29"
$ws.Range("E46").Value = "  -7.34%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "142.54"
$ws.Range("E47").Value = "  -1.97%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.528"
$ws.Range("E48").Value = "  -0.89%  "
$ws.Range("E49").Value = "  -2.13%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.58"
$ws.Range("E50").Value = "  -3.50%  "
$ws.Range("E51").Value = "  -1.07%  "
